$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report record needs to be inserted as row 182. Insert a blank
# row there (shifting the existing rows 182:227 down to 183:228, growing the
# used range from A1:R227 to A1:R228) and then fill in the new record.
$ws.Rows("182:182").Insert(-4121)  # -4121 = xlShiftDown

$ws.Range("D182").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = 44551
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = 100112045
$ws.Range("G182").Value = "Zapallo"
$ws.Range("H182").Value = "Paine"
$ws.Range("I182").Value = "1a nueva(o)"
$ws.Range("J182").Value = 1200
$ws.Range("K182").Value = 500
$ws.Range("L182").Value = 500
$ws.Range("M182").Value = 500
$ws.Range("N182").Value = "$/kilo (volumen en unidades)"
$ws.Range("O182").Value = "Región de O'Higgins"
$ws.Range("P182").Value = 500
$ws.Range("Q182").Value = 1
$ws.Range("R182").Value = "Hortaliza"
